$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the remaining "door" relationship opinions (S3, Y2:Y4, AB2:AB4)
# Order matches the order these strings were first authored so that new
# shared-string entries line up with the source workbook.
$ws.Range("S3").Value = "Takes muscular and silent nature as a trheat."

$ws.Range("Y2").Value = "Wants to let in if Bob isn't present, otherwise will want to only let them crash for the night."
$ws.Range("AB2").Value = "Wants to let in if Bob isn't present, otherwise will want to only let them crash for the night."

$ws.Range("Y3").Value = "If Jessica is present he will have the same opinions. Otherwise, he will be against letting them in - too many mouths to feed."
$ws.Range("AB3").Value = "If Jessica is present he will have the same opinions. Otherwise, he will be against letting them in - too many mouths to feed."

$ws.Range("Y4").Value = "Violet will be in favour if they have the car - more skilled members to join the party. Otherwise will have similar opinion to Bob."
$ws.Range("AB4").Value = "Violet will be in favour if they have the car - more skilled members to join the party. Otherwise will have similar opinion to Bob."

# Update the saved view position/selection to match where the author left off
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 21
$ws.Range("AA4").Select()
